$d = $word.ActiveDocument

# Locate the end of the "TZU..." bibliography paragraph (the last one we keep).
$tzuRange = $d.Content.Duplicate
$foundTzu = $tzuRange.Find.Execute(
    "TZU, S. A Arte da Guerra (Edição Completa). São Paulo. WMF Martins Fontes, 2009.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $foundTzu) {
    throw "Could not find the TZU bibliography paragraph."
}

# Locate the end of the trailing copyright/footer paragraph (the last one we remove).
$copyrightRange = $d.Content.Duplicate
$foundCopyright = $copyrightRange.Find.Execute(
    "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $foundCopyright) {
    throw "Could not find the trailing copyright paragraph."
}

# Build a range that covers everything between the two (inclusive of their
# own paragraph marks), i.e. the blank paragraph, the "Ver no Jupiter..."
# paragraph and the "© 2020..." paragraph, then delete it in one go.
$deleteStart = $tzuRange.End + 1
$deleteEnd = $copyrightRange.End + 1

$toDelete = $d.Range($deleteStart, $deleteEnd)
$toDelete.Delete()
